# Actor_RobotObj_01_Ksetting: add left-arm IK setting rows (mirrors the
# existing right-arm block in rows 4-9 down into rows 15-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the cell formatting (styles/borders) from the already-populated
#    "right arm" block (A4:N9) down onto the currently-blank A15:N20 block,
#    so the new rows pick up the same look (this matches the style indices
#    seen in the target file for rows 15-20).
$ws.Range("A4:N9").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# A18's first cell keeps the source block's row-7 border variant after the
# block paste (row 7 is the first row of its visual sub-group and uses a
# slightly different top border), but row 18 should look like the other
# "kp" rows in the new block, so re-stamp A18's format from A15.
$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Fill in the values for the new "left arm" rows.

# Row 15
$ws.Range("A15").Value = 2401
$ws.Range("B15").Value = "LATA01"
$ws.Range("C15").Value = 400
$ws.Range("D15").Value = "LA01"
$ws.Range("E15").Value = 1400
$ws.Range("F15").Value = "LAEE01"
$ws.Range("G15").Value = 2400
$ws.Range("H15").Value = "LATA01"
$ws.Range("I15").Value = "11: diry_look_pos"
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = "on"
$ws.Range("L15").Value = "off"
$ws.Range("M15").Value = "on"
$ws.Range("N15").Value = "Left wrist Rot IK"

# Row 16
$ws.Range("A16").Value = 2402
$ws.Range("B16").Value = "LATA02"
$ws.Range("C16").Value = 401
$ws.Range("D16").Value = "LA02"
$ws.Range("E16").Value = 1400
$ws.Range("F16").Value = "LAEE01"
$ws.Range("G16").Value = 2400
$ws.Range("H16").Value = "LATA01"
$ws.Range("I16").Value = "11: diry_look_pos"
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = "on"
$ws.Range("L16").Value = "off"
$ws.Range("M16").Value = "on"
$ws.Range("N16").Value = "Left wrist Rot IK"

# Row 17
$ws.Range("A17").Value = 2403
$ws.Range("B17").Value = "LATA03"
$ws.Range("C17").Value = 402
$ws.Range("D17").Value = "LA03"
$ws.Range("E17").Value = 1400
$ws.Range("F17").Value = "LAEE01"
$ws.Range("G17").Value = 2400
$ws.Range("H17").Value = "LATA01"
$ws.Range("I17").Value = "11: diry_look_pos"
$ws.Range("J17").Value = 0.2
$ws.Range("K17").Value = "on"
$ws.Range("L17").Value = "off"
$ws.Range("M17").Value = "on"
$ws.Range("N17").Value = "Left wrist Rot IK"

# Row 18
$ws.Range("A18").Value = 2404
$ws.Range("B18").Value = "LATA04"
$ws.Range("C18").Value = 400
$ws.Range("D18").Value = "LA01"
$ws.Range("E18").Value = 1401
$ws.Range("F18").Value = "LAEE02"
$ws.Range("G18").Value = 2401
$ws.Range("H18").Value = "LATA02"
$ws.Range("I18").Value = "0: pos_to_pos"
$ws.Range("J18").Value = 0.2
$ws.Range("K18").Value = "on"
$ws.Range("L18").Value = "off"
$ws.Range("M18").Value = "on"
$ws.Range("N18").Value = "Left elbow IK"

# Row 19
$ws.Range("A19").Value = 2405
$ws.Range("B19").Value = "LATA05"
$ws.Range("C19").Value = 401
$ws.Range("D19").Value = "LA02"
$ws.Range("E19").Value = 1402
$ws.Range("F19").Value = "LAEE03"
$ws.Range("G19").Value = 2402
$ws.Range("H19").Value = "LATA03"
$ws.Range("I19").Value = "0: pos_to_pos"
$ws.Range("J19").Value = 0.8
$ws.Range("K19").Value = "on"
$ws.Range("L19").Value = "off"
$ws.Range("M19").Value = "on"
$ws.Range("N19").Value = "Left wrist Pos IK"

# Row 20
$ws.Range("A20").Value = 2406
$ws.Range("B20").Value = "LATA06"
$ws.Range("C20").Value = 400
$ws.Range("D20").Value = "LA01"
$ws.Range("E20").Value = 1402
$ws.Range("F20").Value = "LAEE03"
$ws.Range("G20").Value = 2402
$ws.Range("H20").Value = "LATA03"
$ws.Range("I20").Value = "0: pos_to_pos"
$ws.Range("J20").Value = 0.1
$ws.Range("K20").Value = "on"
$ws.Range("L20").Value = "on"
$ws.Range("M20").Value = "on"
$ws.Range("N20").Value = "Left wrist Pos IK"

# 3) Leave the selection where the author left it when they saved.
$ws.Range("I13").Select()
